$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '64.360.45'
$ws.Range("E2").Value = '  +2.32%  '
$ws.Range("D3").Value = '3.080.72'
$ws.Range("E3").Value = '  +1.38%  '
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '559.35'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.94%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '145.59'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +6.00%  '
$ws.Range("E7").Value = '  +0.09%  '
$ws.Range("D8").Value = '3.079.04'
$ws.Range("E8").Value = '  +1.44%  '
$ws.Range("E9").Value = '  +1.26%  '
$ws.Range("E10").Value = '  +2.91%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.18'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -1.63%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.471'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +4.82%  '
$ws.Range("E13").Value = '  +1.17%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '35.19'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +2.47%  '
$ws.Range("D15").Value = '3.577.94'
$ws.Range("E15").Value = '  +1.23%  '
$ws.Range("D16").Value = '64.402.97'
$ws.Range("E16").Value = '  +2.33%  '
$ws.Range("D17").Value = '3.077.71'
$ws.Range("E17").Value = '  +1.42%  '
$ws.Range("E18").Value = '  +1.74%  '
$ws.Range("E19").Value = '  +1.31%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '478.42'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.02%  '
$ws.Range("E21").Value = '  +2.64%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.675'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +1.10%  '
$ws.Range("E23").Value = '  +5.54%  '
$ws.Range("E24").Value = '  +10.43%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '81.41'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.54%  '
$ws.Range("E26").Value = '  -0.05%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.80'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +2.35%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.06'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +3.08%  '
$ws.Range("E29").Value = '  +5.33%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.00'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.15%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '26.18'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +1.58%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.14'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.15%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '2.50'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +4.39%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.58'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -1.30%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '55.96'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +1.38%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '6.19'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +4.83%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '459.13'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.07%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.01'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +19.11%  '
$ws.Range("E39").Value = '  +2.43%  '
$ws.Range("E40").Value = '  +3.73%  '
$ws.Range("D41").Value = '2.965.92'
$ws.Range("E41").Value = '  -3.11%  '
$ws.Range("E42").Value = '  +0.56%  '
$ws.Range("E43").Value = '  -1.06%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '27.82'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.05%  '
$ws.Range("E45").Value = '  +4.78%  '
$ws.Range("E46").Value = '  +6.09%  '
$ws.Range("E47").Value = '  +0.02%  '
$ws.Range("E48").Value = '  +2.95%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '121.33'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +4.33%  '
$ws.Range("E50").Value = '  +2.68%  '
$ws.Range("E51").Value = '  +1.52%  '
